$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A174").Value = "IMX-USD"
$ws.Range("A175").Value = "TAO-USD"
$ws.Range("A176").Value = "GRT-USD"
